$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows with D (price) and/or E (volume) text updates ---
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.491.76"
$ws.Range("E2").Value = "  +1.77%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.671.32"
$ws.Range("E3").Value = "  +1.59%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.65"
$ws.Range("E5").Value = "  +2.03%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5279"
$ws.Range("E6").Value = "  +1.11%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2677"
$ws.Range("E8").Value = "  +2.62%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06374"
$ws.Range("E9").Value = "  +0.34%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.74"
$ws.Range("E10").Value = "  +4.55%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07807"
$ws.Range("E11").Value = "  +2.01%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.694.85"
$ws.Range("E12").Value = "  +3.13%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.483"
$ws.Range("E13").Value = "  +1.45%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5561"
$ws.Range("E14").Value = "  +0.45%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅8297"
$ws.Range("E15").Value = "  -0.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.57"
$ws.Range("E16").Value = "  +1.14%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.489.95"
$ws.Range("E17").Value = "  +1.76%  "

# Row 18
$ws.Range("E18").Value = "  +0.07%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.758"
$ws.Range("E19").Value = "  +0.97%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.95"
$ws.Range("E20").Value = "  +2.59%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.31"
$ws.Range("E21").Value = "  +1.34%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.299"
$ws.Range("E22").Value = "  +0.75%  "

# Row 23
$ws.Range("E23").Value = "  +0.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1270"
$ws.Range("E24").Value = "  +4.40%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "137.95"
$ws.Range("E25").Value = "  -5.07%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.377"
$ws.Range("E26").Value = "  -0.31%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.34"
$ws.Range("E27").Value = "  +3.37%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.429"
$ws.Range("E28").Value = "  +2.79%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06238"
$ws.Range("E29").Value = "  +4.89%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.287"
$ws.Range("E30").Value = "  +1.66%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.608"
$ws.Range("E31").Value = "  +5.93%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.416"
$ws.Range("E32").Value = "  +0.59%  "

# Row 33
$ws.Range("E33").Value = "  +2.04%  "

# Row 34
$ws.Range("E34").Value = "  +1.23%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6141"
$ws.Range("E35").Value = "  +9.31%  "

# Row 37
$ws.Range("E37").Value = "  +1.16%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01615"
$ws.Range("E38").Value = "  +0.51%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.037"
$ws.Range("E39").Value = "  +3.15%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.090.99"
$ws.Range("E40").Value = "  +6.45%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8583"
$ws.Range("E41").Value = "  +0.42%  "

# Row 42
$ws.Range("E42").Value = "  -0.01%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.63"
$ws.Range("E43").Value = "  +2.16%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.815.52"
$ws.Range("E44").Value = "  +1.15%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "58.33"
$ws.Range("E45").Value = "  +4.75%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈108"
$ws.Range("E46").Value = "  -2.46%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("E48").Value = "  -0.44%  "

# Row 50
$ws.Range("E50").Value = "  +0.88%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4232"
$ws.Range("E51").Value = "  +0.35%  "

# --- Row 47 / 49 swap (EnergySwap <-> RenderToken) ---
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.521"
$ws.Range("E47").Value = "  +10.20%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.091"
$ws.Range("E49").Value = "  +0.32%  "
